# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
#
# A handful of match rows had their data (id, teams, odds, …) entered
# against the wrong row number. Swap the B:AC payload between the affected
# row pairs so each record lands on its correct row; column A (the running
# index) is left untouched since it belongs to the row, not the record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $row1, $row2) {
    $addr1 = "B$row1" + ":" + "AC$row1"
    $addr2 = "B$row2" + ":" + "AC$row2"

    $rng1 = $ws.Range($addr1)
    $rng2 = $ws.Range($addr2)

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}

Swap-RowData $ws 148 150
Swap-RowData $ws 153 154
